$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 888812
$ws.Range("A3").Value = 888812

$ws.Range("C3").Select()
